$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vehicle identifiers (row 2 and row 3 of the "Patente"/"Motor"/"Chasis" block)
$ws.Range("H2").Value = "ZZZ522"
$ws.Range("I2").Value = "ABC0987AX318"
$ws.Range("J2").Value = "MMAA09XFGS313"

$ws.Range("H3").Value = "ZZZ523"
$ws.Range("I3").Value = "ABC0987AX319"
$ws.Range("J3").Value = "MMAA09XFGS314"

# Move the "Movilidad" / "Hasta $150.000" accessory row from row 2 to row 4
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("L4").Value = "Movilidad"
$ws.Range("M4").Value = "Hasta $150.000"

# Clear the leftover cells on H4:J4 and H5:J5 (content + formatting)
$ws.Range("H4:J5").Clear()

# Update the duplicate-values conditional formatting range (exclude rows 4-5, keep rows 6-18 and 2-3)
$cf = $ws.Range("H2:J18").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("H6:J18,H2:J3"))

# Move the active selection to I7
$ws.Range("I7").Select()
